$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before writing values so that
# numeric-looking strings (e.g. "4.280") are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.162.38'
$ws.Range("E2").Value = '  -4.46%  '
$ws.Range("D3").Value = '1.654.68'
$ws.Range("E3").Value = '  -3.37%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '215.22'
$ws.Range("E5").Value = '  -4.32%  '
$ws.Range("D6").Value = '0.5117'
$ws.Range("E6").Value = '  -3.52%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.2584'
$ws.Range("E8").Value = '  -3.16%  '
$ws.Range("D9").Value = '0.06402'
$ws.Range("E9").Value = '  -4.16%  '
$ws.Range("E10").Value = '  -4.55%  '
$ws.Range("D11").Value = '0.07817'
$ws.Range("E11").Value = '  +1.55%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.658.73'
$ws.Range("E12").Value = '  -3.22%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.280'
$ws.Range("E13").Value = '  -5.28%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.881.87'
$ws.Range("E14").Value = '  -3.41%  '
$ws.Range("D15").Value = '0.5516'
$ws.Range("E15").Value = '  -5.46%  '
$ws.Range("D16").Value = '0.0₅8010'
$ws.Range("E16").Value = '  -2.47%  '
$ws.Range("D17").Value = '63.84'
$ws.Range("E17").Value = '  -6.25%  '
$ws.Range("D18").Value = '26.181.73'
$ws.Range("E18").Value = '  -4.46%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '208.54'
$ws.Range("E20").Value = '  -5.76%  '
$ws.Range("D21").Value = '4.412'
$ws.Range("E21").Value = '  -5.00%  '
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").Value = '6.015'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '143.11'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  +2.39%  '
$ws.Range("D27").Value = '0.1165'
$ws.Range("E27").Value = '  -3.92%  '
$ws.Range("D28").Value = '6.968'
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("E29").Value = '  -3.03%  '
$ws.Range("D30").Value = '0.05127'
$ws.Range("E30").Value = '  -4.36%  '
$ws.Range("D31").Value = '1.241'
$ws.Range("E31").Value = '  -4.48%  '
$ws.Range("D32").Value = '3.342'
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").Value = '3.216'
$ws.Range("E33").Value = '  -6.39%  '
$ws.Range("E34").Value = '  -5.05%  '
$ws.Range("E35").Value = '  -4.20%  '
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").Value = '0.9281'
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("D38").Value = '0.5679'
$ws.Range("E38").Value = '  -3.22%  '
$ws.Range("D39").Value = '1.155.66'
$ws.Range("E39").Value = '  +5.96%  '
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '0.8333'
$ws.Range("E42").Value = '  -1.64%  '
$ws.Range("D43").Value = '5.641'
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("D44").Value = '100.15'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = '1.791.74'
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").Value = '55.67'
$ws.Range("E48").Value = '  -3.95%  '
$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '7.846'
$ws.Range("E50").Value = '  -3.01%  '
$ws.Range("D51").Value = '0.05037'
$ws.Range("E51").Value = '  -3.89%  '

# Remove the temporary text-format styling so the cells are left without
# an explicit style, matching the original (unstyled) cells while keeping
# their content as text.
$priceRange.ClearFormats()

Write-Host "Crypto prices updated"
